$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 572.1
$ws.Range("I41").Value = 388.18182
$ws.Range("J41").Value = 796.8889
$ws.Range("K41").Value = 388.18182
$ws.Range("L41").Value = 796.8889
$ws.Range("M41").Value = 51.81817999999998
$ws.Range("N41").Value = -1676.8889

# Row 111
$ws.Range("H111").Value = 18750
$ws.Range("I111").Value = 15000
$ws.Range("J111").Value = 30000
$ws.Range("K111").Value = 45000
$ws.Range("L111").Value = 90000
$ws.Range("M111").Value = -41933
$ws.Range("N111").Value = -96134

# Row 125
$ws.Range("H125").Value = 3895.5
$ws.Range("I125").Value = 1221.3334
$ws.Range("J125").Value = 5500
$ws.Range("K125").Value = 10992.0006
$ws.Range("L125").Value = 49500
$ws.Range("M125").Value = -8532.000599999999
$ws.Range("N125").Value = -54420

# Row 135
$ws.Range("H135").Value = 2183.724
$ws.Range("I135").Value = 419.92593
$ws.Range("J135").Value = 25995
$ws.Range("K135").Value = 3779.33337
$ws.Range("L135").Value = 233955
$ws.Range("M135").Value = -1244.33337
$ws.Range("N135").Value = -239025

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6112.625
$ws.Range("I45").Value = 6118.1113
$ws.Range("J45").Value = 6105.5713
$ws.Range("K45").Value = 6118.1113
$ws.Range("L45").Value = 6105.5713
$ws.Range("M45").Value = -5741.1113
$ws.Range("N45").Value = -6859.5713

# Row 61
$ws.Range("H61").Value = 2725.4412
$ws.Range("I61").Value = 1763.7778
$ws.Range("J61").Value = 3807.3125
$ws.Range("K61").Value = 1763.7778
$ws.Range("L61").Value = 3807.3125
$ws.Range("M61").Value = -1551.7778
$ws.Range("N61").Value = -4231.3125

# Row 88
$ws.Range("H88").Value = 1957.375
$ws.Range("I88").Value = 1861.5
$ws.Range("J88").Value = 2245
$ws.Range("K88").Value = 1861.5
$ws.Range("L88").Value = 2245
$ws.Range("M88").Value = -1455.5
$ws.Range("N88").Value = -3057

# Row 91
$ws.Range("H91").Value = 1957.375
$ws.Range("I91").Value = 1861.5
$ws.Range("J91").Value = 2245
$ws.Range("K91").Value = 1861.5
$ws.Range("L91").Value = 2245
$ws.Range("M91").Value = -457.5
$ws.Range("N91").Value = -5053

# Row 110
$ws.Range("H110").Value = 3535.7856
$ws.Range("I110").Value = 1662.625
$ws.Range("J110").Value = 6033.3335
$ws.Range("K110").Value = 1662.625
$ws.Range("L110").Value = 6033.3335
$ws.Range("M110").Value = 382.375
$ws.Range("N110").Value = -10123.3335

# Row 122
$ws.Range("H122").Value = 2364.08
$ws.Range("I122").Value = 1662
$ws.Range("K122").Value = 4986
$ws.Range("M122").Value = -2536

# Row 136
$ws.Range("H136").Value = 2725.4412
$ws.Range("I136").Value = 1763.7778
$ws.Range("J136").Value = 3807.3125
$ws.Range("K136").Value = 5291.3334
$ws.Range("L136").Value = 11421.9375
$ws.Range("M136").Value = -2741.3334
$ws.Range("N136").Value = -16521.9375

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1076.5
$ws.Range("I94").Value = 775.3333
$ws.Range("J94").Value = 1980
$ws.Range("K94").Value = 775.3333
$ws.Range("L94").Value = 1980
$ws.Range("M94").Value = -324.3333
$ws.Range("N94").Value = -2882

# Row 107
$ws.Range("H107").Value = 1726.258
$ws.Range("I107").Value = 1648.56
$ws.Range("J107").Value = 2050
$ws.Range("K107").Value = 1648.56
$ws.Range("L107").Value = 2050
$ws.Range("M107").Value = 271.4400000000001
$ws.Range("N107").Value = -5890

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2839.5386
$ws.Range("I31").Value = 2281.158
$ws.Range("J31").Value = 3370
$ws.Range("K31").Value = 2281.158
$ws.Range("L31").Value = 3370
$ws.Range("M31").Value = -1986.158
$ws.Range("N31").Value = -3960

# Row 34
$ws.Range("H34").Value = 2839.5386
$ws.Range("I34").Value = 2281.158
$ws.Range("J34").Value = 3370
$ws.Range("K34").Value = 2281.158
$ws.Range("L34").Value = 3370
$ws.Range("M34").Value = -2079.158
$ws.Range("N34").Value = -3774

# Row 58
$ws.Range("H58").Value = 1474945.4
$ws.Range("I58").Value = 3690.5334
$ws.Range("J58").Value = 2636462.2
$ws.Range("K58").Value = 3690.5334
$ws.Range("L58").Value = 2636462.2
$ws.Range("M58").Value = -3487.5334
$ws.Range("N58").Value = -2636868.2

# Row 136
$ws.Range("H136").Value = 1474945.4
$ws.Range("I136").Value = 3690.5334
$ws.Range("J136").Value = 2636462.2
$ws.Range("K136").Value = 11071.6002
$ws.Range("L136").Value = 7909386.600000001
$ws.Range("M136").Value = -8521.600199999999
$ws.Range("N136").Value = -7914486.600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 2500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 2500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 7500
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -7948

# Row 132
$ws.Range("H132").Value = 2084.7
$ws.Range("J132").Value = 2288.125
$ws.Range("L132").Value = 20593.125
$ws.Range("N132").Value = -25653.125

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2978414.2
$ws.Range("I102").Value = 5103225.5
$ws.Range("J102").Value = 3678.2
$ws.Range("K102").Value = 5103225.5
$ws.Range("L102").Value = 3678.2
$ws.Range("M102").Value = -5101603.5
$ws.Range("N102").Value = -6922.2

# Row 132
$ws.Range("H132").Value = 2715.3333
$ws.Range("I132").Value = 3120.8
$ws.Range("J132").Value = 2512.6
$ws.Range("K132").Value = 9362.400000000001
$ws.Range("L132").Value = 7537.799999999999
$ws.Range("M132").Value = -6832.400000000001
$ws.Range("N132").Value = -12597.8

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3016.6667
$ws.Range("J81").Value = 1524.75
$ws.Range("L81").Value = 3049.5
$ws.Range("N81").Value = -5171.5

# Row 84
$ws.Range("H84").Value = 3016.6667
$ws.Range("J84").Value = 1524.75
$ws.Range("L84").Value = 15247.5
$ws.Range("N84").Value = -25855.5

# Row 136
$ws.Range("H136").Value = 13890938
$ws.Range("I136").Value = 22728416
$ws.Range("J136").Value = 3471.4285
$ws.Range("K136").Value = 68185248
$ws.Range("L136").Value = 10414.2855
$ws.Range("M136").Value = -68182698
$ws.Range("N136").Value = -15514.2855

# Row 138
$ws.Range("H138").Value = 47429
$ws.Range("J138").Value = 47429
$ws.Range("L138").Value = 47429
$ws.Range("N138").Value = -57709
